$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 235, shifting existing rows 235:280 down to 236:281
$ws.Rows.Item(235).Insert()

# Populate the newly inserted row 235 with its data.
# Columns A,B,C,E,F,G,H,I,N,O,Q,R keep the same values the old row 235 had
# (same market/category/quality/unit/origin), while D,J,K,L,M,P hold the
# new values introduced by this edit.
$ws.Range("A235").Value = 10
$ws.Range("B235").Value = "Vega Modelo de Temuco"
$ws.Range("C235").Value = "La Araucanía"
$ws.Range("D235").Value = 44694
$ws.Range("D235").NumberFormat = $ws.Range("D236").NumberFormat
$ws.Range("E235").Value = 9
$ws.Range("F235").Value = 100112001
$ws.Range("G235").Value = "Berenjena"
$ws.Range("H235").Value = "Sin especificar"
$ws.Range("I235").Value = "Primera"
$ws.Range("J235").Value = 45
$ws.Range("K235").Value = 9000
$ws.Range("L235").Value = 10000
$ws.Range("M235").Value = 9444
$ws.Range("N235").Value = "`$/caja 60 unidades"
$ws.Range("O235").Value = "Región de Arica y Parinacota"
$ws.Range("P235").Value = 157
$ws.Range("Q235").Value = 60
$ws.Range("R235").Value = "Hortaliza"
